$d = $word.ActiveDocument

# --- Merge the split runs in the Title, Author and Abstract paragraphs into
# --- a single run each (their joined text is unchanged; only the run
# --- splitting collapses) by doing an exact-text Find & Replace scoped to
# --- each paragraph's own Range, so identical text elsewhere (e.g. the
# --- "Sophie Chowgule" that also appears in the version-history note) is
# --- left untouched. Word's Find/Replace rebuilds the matched range as a
# --- single run.

$titlePar = $d.Paragraphs(1).Range
$titlePar.Find.Execute(
    "Questions: PMFs, PDFs, and CDFs", $true, $false, $false, $false, $false,
    $true, 1, $false, "Questions: PMFs, PDFs, and CDFs", 2)

$authorPar = $d.Paragraphs(2).Range
$authorPar.Find.Execute(
    "Sophie Chowgule", $true, $false, $false, $false, $false,
    $true, 1, $false, "Sophie Chowgule", 2)

$abstractPar = $d.Paragraphs(4).Range
$abstractPar.Find.Execute(
    "A selection of questions to test your understanding of Probability Mass Functions (PMFs), Probability Density Functions (PDFs), and Cumulative Distribution Functions (CDFs).",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "A selection of questions to test your understanding of Probability Mass Functions (PMFs), Probability Density Functions (PDFs), and Cumulative Distribution Functions (CDFs).",
    2)

# --- Drop the redundant explicit "left" paragraph alignment on every
# --- "Compact" styled paragraph (table-cell math/labels). Left is already
# --- the inherited default, so (re)setting it clears the direct <w:jc/>.

foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Compact") {
        $p.Alignment = 0
    }
}
